$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Save" header in H1, copying the header style/format from G1
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Fill the new "Save" column values for the two data rows
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1

$excel.CutCopyMode = 0
